$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "58.936.28"
Set-TextValue 2 5 "  +2.56%  "

Set-TextValue 3 4 "2.494.40"
Set-TextValue 3 5 "  +2.73%  "

Set-TextValue 4 4 "1.00"
Set-TextValue 4 5 "  +0.06%  "

Set-TextValue 5 4 "533.97"
Set-TextValue 5 5 "  +4.50%  "

Set-TextValue 6 4 "135.37"
Set-TextValue 6 5 "  +5.02%  "

Set-TextValue 7 5 "  -0.09%  "

Set-TextValue 8 4 "0.566"
Set-TextValue 8 5 "  +3.31%  "

Set-TextValue 9 4 "2.518.65"
Set-TextValue 9 5 "  +3.29%  "

Set-TextValue 10 5 "  +5.22%  "

Set-TextValue 11 5 "  -1.40%  "

Set-TextValue 12 4 "5.25"
Set-TextValue 12 5 "  +2.12%  "

Set-TextValue 13 5 "  +1.72%  "

Set-TextValue 14 4 "2.941.15"
Set-TextValue 14 5 "  +2.87%  "

Set-TextValue 15 4 "58.913.89"
Set-TextValue 15 5 "  +2.68%  "

Set-TextValue 16 4 "22.48"
Set-TextValue 16 5 "  +3.71%  "

Set-TextValue 17 5 "  +3.89%  "

Set-TextValue 18 4 "2.507.44"
Set-TextValue 18 5 "  +3.01%  "

Set-TextValue 19 4 "10.71"
Set-TextValue 19 5 "  +3.07%  "

Set-TextValue 20 5 "  +3.92%  "

Set-TextValue 21 4 "322.17"
Set-TextValue 21 5 "  +2.43%  "

Set-TextValue 22 4 "6.16"
Set-TextValue 22 5 "  +9.65%  "

Set-TextValue 23 5 "  -0.36%  "

Set-TextValue 24 4 "65.68"
Set-TextValue 24 5 "  +3.87%  "

Set-TextValue 25 5 "  +2.13%  "

Set-TextValue 26 4 "0.996"
Set-TextValue 26 5 "  -0.34%  "

Set-TextValue 27 4 "0.160"
Set-TextValue 27 5 "  +1.41%  "

Set-TextValue 28 4 "7.51"
Set-TextValue 28 5 "  +4.31%  "

Set-TextValue 29 4 "0.0₃0768"
Set-TextValue 29 5 "  +7.24%  "

Set-TextValue 30 4 "171.75"
Set-TextValue 30 5 "  +1.11%  "

Set-TextValue 31 5 "  +5.19%  "

Set-TextValue 32 5 "  +4.76%  "

Set-TextValue 33 4 "6.31"
Set-TextValue 33 5 "  +1.79%  "

Set-TextValue 34 5 "  +0.03%  "

Set-TextValue 35 4 "0.994"
Set-TextValue 35 5 "  -0.37%  "

Set-TextValue 36 4 "18.22"
Set-TextValue 36 5 "  +3.44%  "

Set-TextValue 37 5 "  -1.07%  "

Set-TextValue 38 4 "4.00"
Set-TextValue 38 5 "  +2.39%  "

Set-TextValue 40 4 "36.70"
Set-TextValue 40 5 "  +1.34%  "

Set-TextValue 41 4 "0.789"
Set-TextValue 41 5 "  +3.04%  "

Set-TextValue 44 4 "3.50"
Set-TextValue 44 5 "  +4.13%  "

Set-TextValue 45 4 "132.09"
Set-TextValue 45 5 "  +10.28%  "

Set-TextValue 46 5 "  +2.99%  "

Set-TextValue 47 5 "  +3.06%  "

Set-TextValue 48 5 "  +6.02%  "

Set-TextValue 49 5 "  +4.90%  "

Set-TextValue 50 4 "17.17"
Set-TextValue 50 5 "  +4.44%  "

Set-TextValue 51 4 "1.761.39"
Set-TextValue 51 5 "  +3.54%  "

# Rows 42 and 43 swap coin/link/price/volume (RenderToken and Bittensor switch ranking order)
Set-TextValue 42 2 "RenderToken"
Set-TextValue 42 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 42 4 "5.22"
Set-TextValue 42 5 "  +7.37%  "

Set-TextValue 43 2 "Bittensor"
Set-TextValue 43 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 43 4 "280.69"
Set-TextValue 43 5 "  +3.55%  "

